$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.955.97'
$ws.Range('E2').Value = '  +1.60%  '
$ws.Range('D3').Value = '1.640.32'
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.44'
$ws.Range('E5').Value = '  -0.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.525'
$ws.Range('E6').Value = '  -1.48%  '
$ws.Range('E7').Value = '  -0.33%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.32'
$ws.Range('E8').Value = '  +0.69%  '
$ws.Range('E9').Value = '  +2.58%  '
$ws.Range('E10').Value = '  +0.33%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0894'
$ws.Range('E11').Value = '  +0.71%  '
$ws.Range('D12').Value = '1.871.92'
$ws.Range('E12').Value = '  -0.13%  '
$ws.Range('D13').Value = '1.635.91'
$ws.Range('E13').Value = '  -0.42%  '
$ws.Range('E14').Value = '  +0.94%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.561'
$ws.Range('E15').Value = '  -3.40%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.62'
$ws.Range('E16').Value = '  +0.66%  '
$ws.Range('D17').Value = '27.906.44'
$ws.Range('E17').Value = '  +1.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '233.48'
$ws.Range('E18').Value = '  +1.83%  '
$ws.Range('E19').Value = '  +0.11%  '
$ws.Range('E20').Value = '  +1.17%  '
$ws.Range('E22').Value = '  +0.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.99'
$ws.Range('E23').Value = '  +3.01%  '
$ws.Range('E24').Value = '  +4.44%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '150.47'
$ws.Range('E25').Value = '  +0.98%  '
$ws.Range('E26').Value = '  -0.38%  '
$ws.Range('E27').Value = '  -0.72%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.69'
$ws.Range('E28').Value = '  +1.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('E30').Value = '  +0.18%  '
$ws.Range('E31').Value = '  -0.43%  '
$ws.Range('E32').Value = '  +0.55%  '
$ws.Range('D33').Value = '1.473.05'
$ws.Range('E33').Value = '  +3.99%  '
$ws.Range('E34').Value = '  -1.68%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.55'
$ws.Range('E35').Value = '  -2.30%  '
$ws.Range('E36').Value = '  -0.55%  '
$ws.Range('E37').Value = '  -0.45%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.880'
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.926'
$ws.Range('E39').Value = '  +13.68%  '
$ws.Range('E40').Value = '  +0.51%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '69.08'
$ws.Range('E41').Value = '  +6.59%  '
$ws.Range('E42').Value = '  -0.25%  '
$ws.Range('E43').Value = '  -2.03%  '
$ws.Range('E44').Value = '  -0.21%  '
$ws.Range('E45').Value = '  +0.25%  '
$ws.Range('E46').Value = '  -0.27%  '
$ws.Range('D47').Value = '1.781.55'
$ws.Range('E47').Value = '  -0.17%  '
$ws.Range('E48').Value = '  +2.65%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '87.26'
$ws.Range('E49').Value = '  +1.47%  '
$ws.Range('E50').Value = '  -1.06%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.84'
$ws.Range('E51').Value = '  +1.62%  '
